$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 402; existing rows 402..462 shift down to 404..464
$ws.Range("A402:A403").EntireRow.Insert()

# New row 402 data
$ws.Cells.Item(402, 1).Value = 10
$ws.Cells.Item(402, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(402, 3).Value = "La Araucanía"
$ws.Cells.Item(402, 4).Value = 45127
$ws.Cells.Item(402, 5).Value = 9
$ws.Cells.Item(402, 6).Value = "Fruta"
$ws.Cells.Item(402, 7).Value = 100102
$ws.Cells.Item(402, 8).Value = "Cítricos"
$ws.Cells.Item(402, 9).Value = 100102006
$ws.Cells.Item(402, 10).Value = "Pomelo"
$ws.Cells.Item(402, 11).Value = "Start Ruby"
$ws.Cells.Item(402, 12).Value = "Primera"
$ws.Cells.Item(402, 13).Value = 125
$ws.Cells.Item(402, 14).Value = 15000
$ws.Cells.Item(402, 15).Value = 15000
$ws.Cells.Item(402, 16).Value = 15000
$ws.Cells.Item(402, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(402, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(402, 19).Value = 1000
$ws.Cells.Item(402, 20).Value = 15

# New row 403 data
$ws.Cells.Item(403, 1).Value = 10
$ws.Cells.Item(403, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(403, 3).Value = "La Araucanía"
$ws.Cells.Item(403, 4).Value = 45127
$ws.Cells.Item(403, 5).Value = 9
$ws.Cells.Item(403, 6).Value = "Fruta"
$ws.Cells.Item(403, 7).Value = 100102
$ws.Cells.Item(403, 8).Value = "Cítricos"
$ws.Cells.Item(403, 9).Value = 100102006
$ws.Cells.Item(403, 10).Value = "Pomelo"
$ws.Cells.Item(403, 11).Value = "Start Ruby"
$ws.Cells.Item(403, 12).Value = "Primera"
$ws.Cells.Item(403, 13).Value = 55
$ws.Cells.Item(403, 14).Value = 17000
$ws.Cells.Item(403, 15).Value = 17000
$ws.Cells.Item(403, 16).Value = 17000
$ws.Cells.Item(403, 17).Value = "$/caja 14 kilos empedrada"
$ws.Cells.Item(403, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(403, 19).Value = 1214
$ws.Cells.Item(403, 20).Value = 14

# Apply the date number format (same as other Fecha cells) to the new D cells
$ws.Cells.Item(402, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(403, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
